# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new rows 82-83) above the previous
# latest entries (old rows 82-83, which shift down to 84-85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 82; existing rows 82:83 shift to 84:85.
$ws.Range("A82:A83").EntireRow.Insert()

# --- New row 82: Primera ---
$ws.Cells.Item(82, 1).Value = 1
$ws.Cells.Item(82, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(82, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(82, 4).Value = 44610
$ws.Cells.Item(82, 5).Value = 15
$ws.Cells.Item(82, 6).Value = "Fruta"
$ws.Cells.Item(82, 7).Value = 100108
$ws.Cells.Item(82, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(82, 9).Value = 100108003
$ws.Cells.Item(82, 10).Value = "Maracuyá"
$ws.Cells.Item(82, 11).Value = "Sin especificar"
$ws.Cells.Item(82, 12).Value = "Primera"
$ws.Cells.Item(82, 13).Value = 130
$ws.Cells.Item(82, 14).Value = 34000
$ws.Cells.Item(82, 15).Value = 35000
$ws.Cells.Item(82, 16).Value = 34500
$ws.Cells.Item(82, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(82, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(82, 19).Value = 1725
$ws.Cells.Item(82, 20).Value = 20

# --- New row 83: Segunda ---
$ws.Cells.Item(83, 1).Value = 1
$ws.Cells.Item(83, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(83, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(83, 4).Value = 44610
$ws.Cells.Item(83, 5).Value = 15
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100108
$ws.Cells.Item(83, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(83, 9).Value = 100108003
$ws.Cells.Item(83, 10).Value = "Maracuyá"
$ws.Cells.Item(83, 11).Value = "Sin especificar"
$ws.Cells.Item(83, 12).Value = "Segunda"
$ws.Cells.Item(83, 13).Value = 150
$ws.Cells.Item(83, 14).Value = 30000
$ws.Cells.Item(83, 15).Value = 31000
$ws.Cells.Item(83, 16).Value = 30500
$ws.Cells.Item(83, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(83, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(83, 19).Value = 1525
$ws.Cells.Item(83, 20).Value = 20
